# Pokemon Munchkin Rules - "Added more epic stuff"
# Edits the rules text box (shape 3 on slide 1):
#   - rewords two bullet paragraphs (Pokemon stages / once-per-game recharge)
#   - bumps every run to 17pt
#   - repositions/resizes the text box slightly
#
# NOTE: this COM-interop's TextRange.InsertBefore/InsertAfter return a range
# object whose .Start/.Length still reflect the ORIGINAL (pre-insert) anchor,
# not the freshly-inserted text. So after any mutating call (.Text = ..,
# .InsertBefore(..), .InsertAfter(..)) we never reuse that object's position -
# we instead recompute plain integer offsets ourselves and fetch brand new
# ranges via $tr.Characters(pos, len).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# ---------------------------------------------------------------------------
# 1) Paragraph 2 ("Pokémon count as items...")
#    "Pokémon count as items. Stage 3, Pseudo-Legendary, Legendary and
#     Mythical Pokémon count as big items"
#    ->
#    "Stage 1 and Stage 2, Pokémon count as items. Stage 3, Pseudo-Legendary,
#     Legendary, Mythical and Fossil Pokémon count as big items"
# ---------------------------------------------------------------------------

# 1a. Prepend "Stage 1 and Stage 2, " before "Pokémon count as items."
$anchor = $tr.Find("Pokémon count as ")
$pos = $anchor.Start
[void]$anchor.InsertBefore("Stage 1 and Stage 2, ")
$tr.Characters($pos, 7).Font.Underline = $true        # "Stage 1"
$tr.Characters($pos + 12, 7).Font.Underline = $true    # "Stage 2"

# 1b. "Legendary and Mythical" -> "Legendary, Mythical"
$anchor = $tr.Find("Legendary and Mythical")
$pos = $anchor.Start
$tr.Characters($pos + 9, 5).Text = ", "

# 1c. "Mythical Pokémon count as " -> "Mythical and Fossil Pokémon count as "
$anchor = $tr.Find("Mythical")
$pos = $anchor.Start + $anchor.Length
$tr.Characters($pos, 18).Text = " and Fossil Pokémon count as "
$tr.Characters($pos + 5, 6).Font.Underline = $true     # "Fossil"

# ---------------------------------------------------------------------------
# 2) Paragraph 4 ("Once per game abilities can be recharged by ...")
#    "rolling a 1d20 and landing 18+"
#    ->
#    "winning in a battle against a monster with bonuses of 20+"
# ---------------------------------------------------------------------------

$onceAnchor = $tr.Find("Once per game")
$afterOnce = $onceAnchor.Start + $onceAnchor.Length

$anchor = $tr.Find("rolling a 1d20 ", $afterOnce)
$anchor.Text = "winning "

$anchor = $tr.Find("and landing 18+", $afterOnce)
$pos = $anchor.Start
$anchor.Text = "in a battle against a monster with bonuses of 20+"
$tr.Characters($pos + 35, 14).Font.Underline = $true   # "bonuses of 20+"

# ---------------------------------------------------------------------------
# 3) Bump every run in the box to 17pt.
# ---------------------------------------------------------------------------
$tr.Font.Size = 17

# ---------------------------------------------------------------------------
# 4) Reposition / resize the text box.  (EMU -> points: / 12700)
# ---------------------------------------------------------------------------
$sh.Left   = 2905458 / 12700
$sh.Top    = 1976213 / 12700
$sh.Width  = 6368381 / 12700
$sh.Height = 3493264 / 12700
